# questionbank.xlsx update
# - R5 (rows 20-22) now asks about Mycenaean palaces / hero shrine / Queen of Tiryns
# - R6 (rows 23-25) now asks about Panhellenic games/shrines
# - R7 (rows 26-28) now asks about the French wine industry / Cyrene / Pindar
# - R8 (rows 29-31) now asks about Clearchus / Spartan hoplite warfare
# - E2/E3 switched from plain shared-string answers to LOWER() formulas
# - Row 20 grew an extra answer line, so its height increased
# - Window scrolled down to the bottom of the sheet (row 25) with E31 selected

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 / Row 3: correct-answer cells become formulas ---
$ws.Range("E2").Formula = '=LOWER("d")'
$ws.Range("E3").Formula = '=LOWER("c")'

# --- R5 block (rows 20-22) ---
$ws.Range("C20").Value = "1. How did the features of Mycenaean palaces reflect the priorities of the royal family?"
$ws.Range("D20").Value = "A) They had luxurious items indicating wealth and comfort`nB) They were built primarily for defense`nC) They were simple and lacked decoration`nD) They focused on agricultural storage`nE) They were designed for religious ceremonies only"
$ws.Range("E20").Value = "a"

$ws.Range("C21").Value = "2. The hero shrine at Therapne was made for whom?"
$ws.Range("D21").Value = "A) Achilles`nB) Odysseus`nC) Agamemnon`nD) Hector`nE) Helen"
$ws.Range("E21").Value = "e"

$ws.Range("C22").Value = "3. What is interesting about the Queen of Tiryns?"
$ws.Range("D22").Value = "A) Leadership`nB) Wealth`nC) False accusation`nD) Military skills`nE) Diplomacy"
$ws.Range("E22").Value = "c"

# Row 20 now wraps to an extra line in column D
$ws.Rows.Item(20).RowHeight = 86.4

# --- R6 block (rows 23-25) ---
$ws.Range("C23").Value = "1. What was the primary focus of the first competitions at Delphi?"
$ws.Range("D23").Value = "A) Wrestling`nB) Musical`nC) Running`nD) Javelin throwing`nE) Chariot racing"
$ws.Range("E23").Value = "b"

$ws.Range("C24").Value = "2. How did Panhellenic shrines serve Greece?"
$ws.Range("D24").Value = "A) Promoted economic growth`nB) Established military alliances`nC) Managed relationships`nD) Enforced legal systems`nE) Expanded territories"
$ws.Range("E24").Value = "c"

$ws.Range("C25").Value = "3. What did the Greeks do to gather yearly?"
$ws.Range("D25").Value = "A) Trade in markets`nB) Fight in wars`nC) Participate in games`nD) Hold religious festivals`nE) Attend political assemblies"
$ws.Range("E25").Value = "c"

# --- R7 block (rows 26-28) ---
$ws.Range("C26").Value = "1. What does the marriage in the passage with the French wine industry symbolize?"
$ws.Range("D26").Value = "A) Greek conquest`nB) Decline of Ligurians`nC) New Greek dynasty`nD) Union of cultures`nE) Political control"
$ws.Range("E26").Value = "d"

$ws.Range("C27").Value = "2. What values are reflected in Pindara's portrayal of Cyrene?"
$ws.Range("D27").Value = "A) Wealth and luxury`nB) Agricultural skills`nC) Strength and heroism`nD) Peace and coexistence`nE) Democratic principles"
$ws.Range("E27").Value = "c"

$ws.Range("C28").Value = "3. Who fell in love with Cyrene according to Pindar's poem?"
$ws.Range("D28").Value = "A) Zeus`nB) Heracles`nC) Apollo`nD) Poseidon`nE) Ares"
$ws.Range("E28").Value = "c"

# --- R8 block (rows 29-31) ---
$ws.Range("C29").Value = "1. When was Clearchus happy?"
$ws.Range("D29").Value = "A) During peace times`nB) When spending his fortune`nC) When fighting`nD) When forming close relationships`nE) When following orders"
$ws.Range("E29").Value = "c"

$ws.Range("C30").Value = "2. What qualities contributed to Spartan's success in hoplite warfare?"
$ws.Range("D30").Value = "A) Superior weaponry`nB) Large armies`nC) Advanced technology`nD) Morale and organization`nE) Naval dominance"
$ws.Range("E30").Value = "d"

$ws.Range("C31").Value = "3. What would stop Spartans from going to battle?"
$ws.Range("D31").Value = "A) Lack of weapons`nB) Poor weather`nC) Insufficient troops`nD) Religious omens`nE) Political decisions"
$ws.Range("E31").Value = "d"

# --- Scroll/selection state as last saved ---
$ws.Range("E31").Select()
